# Parameter tuning for GA configuration
# Update E (Thu), F (Tiet BD) and H (Diem) columns for rows 2-99
# to the tuned values per the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(2, 6).Value = 5

# Row 3
$ws.Cells.Item(3, 6).Value = 8

# Row 4
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(4, 6).Value = 9

# Row 5
$ws.Cells.Item(5, 5).Value = 4

# Row 6
$ws.Cells.Item(6, 6).Value = 2

# Row 7
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = 9

# Row 8
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 8

# Row 9
$ws.Cells.Item(9, 5).Value = 2

# Row 10
$ws.Cells.Item(10, 6).Value = 7

# Row 11
$ws.Cells.Item(11, 5).Value = 6
$ws.Cells.Item(11, 6).Value = 10
$ws.Cells.Item(11, 8).Value = 0

# Row 12
$ws.Cells.Item(12, 5).Value = 6
$ws.Cells.Item(12, 6).Value = 2

# Row 13
$ws.Cells.Item(13, 5).Value = 7
$ws.Cells.Item(13, 6).Value = 5

# Row 14
$ws.Cells.Item(14, 5).Value = 7
$ws.Cells.Item(14, 6).Value = 8

# Row 15
$ws.Cells.Item(15, 5).Value = 3

# Row 16
$ws.Cells.Item(16, 5).Value = 5
$ws.Cells.Item(16, 6).Value = 10

# Row 17
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 3
$ws.Cells.Item(17, 8).Value = 0

# Row 18
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 5

# Row 19
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = 9

# Row 20
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = 8

# Row 21
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 11
$ws.Cells.Item(21, 8).Value = 0

# Row 22
$ws.Cells.Item(22, 5).Value = 5
$ws.Cells.Item(22, 6).Value = 4

# Row 24
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 10
$ws.Cells.Item(24, 8).Value = 0

# Row 25
$ws.Cells.Item(25, 5).Value = 3

# Row 26
$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 6).Value = 9

# Row 27
$ws.Cells.Item(27, 6).Value = 3

# Row 28
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 6).Value = 8

# Row 29
$ws.Cells.Item(29, 5).Value = 3

# Row 30
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 6).Value = 4

# Row 32
$ws.Cells.Item(32, 5).Value = 7

# Row 33
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 8

# Row 34
$ws.Cells.Item(34, 6).Value = 4

# Row 35
$ws.Cells.Item(35, 5).Value = 6
$ws.Cells.Item(35, 6).Value = 8

# Row 36
$ws.Cells.Item(36, 5).Value = 6
$ws.Cells.Item(36, 6).Value = 1

# Row 37
$ws.Cells.Item(37, 5).Value = 5
$ws.Cells.Item(37, 6).Value = 7

# Row 38
$ws.Cells.Item(38, 5).Value = 2
$ws.Cells.Item(38, 6).Value = 10

# Row 39
$ws.Cells.Item(39, 5).Value = 7
$ws.Cells.Item(39, 6).Value = 3

# Row 40
$ws.Cells.Item(40, 5).Value = 2
$ws.Cells.Item(40, 6).Value = 10

# Row 41
$ws.Cells.Item(41, 5).Value = 5
$ws.Cells.Item(41, 6).Value = 5

# Row 42
$ws.Cells.Item(42, 5).Value = 4
$ws.Cells.Item(42, 6).Value = 6
$ws.Cells.Item(42, 8).Value = -25

# Row 43
$ws.Cells.Item(43, 5).Value = 6
$ws.Cells.Item(43, 6).Value = 9

# Row 44
$ws.Cells.Item(44, 6).Value = 7

# Row 46
$ws.Cells.Item(46, 5).Value = 3
$ws.Cells.Item(46, 6).Value = 11

# Row 47
$ws.Cells.Item(47, 5).Value = 5

# Row 48
$ws.Cells.Item(48, 5).Value = 6
$ws.Cells.Item(48, 6).Value = 6
$ws.Cells.Item(48, 8).Value = -25

# Row 49
$ws.Cells.Item(49, 5).Value = 3
$ws.Cells.Item(49, 6).Value = 9

# Row 50
$ws.Cells.Item(50, 5).Value = 4
$ws.Cells.Item(50, 6).Value = 10

# Row 51
$ws.Cells.Item(51, 5).Value = 5
$ws.Cells.Item(51, 6).Value = 10

# Row 52
$ws.Cells.Item(52, 5).Value = 4
$ws.Cells.Item(52, 6).Value = 8

# Row 53
$ws.Cells.Item(53, 5).Value = 7
$ws.Cells.Item(53, 6).Value = 3

# Row 54
$ws.Cells.Item(54, 5).Value = 2
$ws.Cells.Item(54, 6).Value = 2

# Row 55
$ws.Cells.Item(55, 5).Value = 3
$ws.Cells.Item(55, 6).Value = 2

# Row 56
$ws.Cells.Item(56, 5).Value = 3

# Row 57
$ws.Cells.Item(57, 5).Value = 7

# Row 58
$ws.Cells.Item(58, 5).Value = 3
$ws.Cells.Item(58, 6).Value = 4

# Row 59
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(59, 6).Value = 4

# Row 60
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = 2

# Row 61
$ws.Cells.Item(61, 5).Value = 4
$ws.Cells.Item(61, 6).Value = 5

# Row 62
$ws.Cells.Item(62, 5).Value = 5

# Row 63
$ws.Cells.Item(63, 6).Value = 9

# Row 64
$ws.Cells.Item(64, 5).Value = 3
$ws.Cells.Item(64, 6).Value = 5

# Row 66
$ws.Cells.Item(66, 5).Value = 2

# Row 67
$ws.Cells.Item(67, 6).Value = 11

# Row 68
$ws.Cells.Item(68, 5).Value = 7
$ws.Cells.Item(68, 6).Value = 7

# Row 69
$ws.Cells.Item(69, 5).Value = 3
$ws.Cells.Item(69, 6).Value = 6
$ws.Cells.Item(69, 8).Value = -25

# Row 70
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 2

# Row 71
$ws.Cells.Item(71, 6).Value = 2

# Row 72
$ws.Cells.Item(72, 5).Value = 4
$ws.Cells.Item(72, 6).Value = 11

# Row 73
$ws.Cells.Item(73, 5).Value = 5
$ws.Cells.Item(73, 6).Value = 3

# Row 74
$ws.Cells.Item(74, 5).Value = 3
$ws.Cells.Item(74, 6).Value = 11

# Row 75
$ws.Cells.Item(75, 5).Value = 6
$ws.Cells.Item(75, 6).Value = 5

# Row 76
$ws.Cells.Item(76, 5).Value = 7

# Row 77
$ws.Cells.Item(77, 5).Value = 5
$ws.Cells.Item(77, 6).Value = 10

# Row 78
$ws.Cells.Item(78, 5).Value = 2
$ws.Cells.Item(78, 6).Value = 5

# Row 79
$ws.Cells.Item(79, 6).Value = 4

# Row 81
$ws.Cells.Item(81, 5).Value = 3
$ws.Cells.Item(81, 6).Value = 4

# Row 82
$ws.Cells.Item(82, 6).Value = 10
$ws.Cells.Item(82, 8).Value = 0

# Row 83
$ws.Cells.Item(83, 5).Value = 6
$ws.Cells.Item(83, 6).Value = 9

# Row 84
$ws.Cells.Item(84, 5).Value = 4

# Row 85
$ws.Cells.Item(85, 5).Value = 7
$ws.Cells.Item(85, 6).Value = 8

# Row 86
$ws.Cells.Item(86, 6).Value = 8

# Row 87
$ws.Cells.Item(87, 5).Value = 5
$ws.Cells.Item(87, 6).Value = 2

# Row 88
$ws.Cells.Item(88, 5).Value = 3
$ws.Cells.Item(88, 6).Value = 10

# Row 89
$ws.Cells.Item(89, 5).Value = 5
$ws.Cells.Item(89, 6).Value = 11

# Row 90
$ws.Cells.Item(90, 5).Value = 3
$ws.Cells.Item(90, 6).Value = 2

# Row 91
$ws.Cells.Item(91, 5).Value = 5
$ws.Cells.Item(91, 6).Value = 3

# Row 92
$ws.Cells.Item(92, 5).Value = 6
$ws.Cells.Item(92, 6).Value = 8

# Row 93
$ws.Cells.Item(93, 5).Value = 2
$ws.Cells.Item(93, 6).Value = 8

# Row 94
$ws.Cells.Item(94, 5).Value = 4
$ws.Cells.Item(94, 6).Value = 9

# Row 95
$ws.Cells.Item(95, 5).Value = 2
$ws.Cells.Item(95, 6).Value = 9

# Row 96
$ws.Cells.Item(96, 5).Value = 5
$ws.Cells.Item(96, 6).Value = 10

# Row 97
$ws.Cells.Item(97, 5).Value = 6
$ws.Cells.Item(97, 6).Value = 5

# Row 98
$ws.Cells.Item(98, 5).Value = 7
$ws.Cells.Item(98, 6).Value = 2

# Row 99
$ws.Cells.Item(99, 5).Value = 3

